$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28, column C: add new reviewer-response note about chapter 2 filtering relocation
$ws.Range("C28").Value = "Some intro from chapter 2 has been removed or replaced into chapter 1, with references added in chapter 2. Section 2.7: Filtering Data has been moved into Chapter 3 as subsection 3.2.3 under methods."

# Row 34, column B: append additional sentence about ozone photolysis cross section
$ws.Range("B34").Value = "Chemistry is not my strongest suit, so I do not pick up errors as readily as I should – thanks for pointing these out. It appears that ozone is photolysed by light up to long wavelengths (~1100nm) however the cross section is orders of magnitude lower for wavelengths past 320nm (reference)."

# Row heights adjust to fit the updated wrapped text
$ws.Rows.Item(28).RowHeight = 256.7
$ws.Rows.Item(34).RowHeight = 256.7

$ws.Range("B34").Select()
